# Add a new row (row 2) to the wishlist with a book entry: "prueba" in
# the Titulo column, Autor/Editorial left blank (but present) - mirrors
# the original inlineStr cells the sheet was authored with.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "prueba"

# Force B2/C2 to exist as (empty) text cells rather than being left
# completely absent from the sheet. A leading apostrophe tells Excel to
# store the cell as text even though there's nothing after it; resetting
# the style back to "Normal" afterwards drops the quote-prefix formatting
# flag that the apostrophe trick would otherwise leave behind, so the
# cells end up blank/unstyled just like their counterparts in row 1's
# layout intent.
$ws.Cells.Item(2, 2).Value = "'"
$ws.Cells.Item(2, 2).Style = "Normal"

$ws.Cells.Item(2, 3).Value = "'"
$ws.Cells.Item(2, 3).Style = "Normal"
